# Add a new worksheet "9" as the first sheet in the workbook, containing
# a fresh submission-group roster (Nitzan Butbul's group).
$wb = $excel.ActiveWorkbook

$firstSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($firstSheet)
$newSheet.Name = "9"

$newSheet.Range("A1").Value = "Nitzan Butbul"
$newSheet.Range("A2").Value = "Hadar Dahan"
$newSheet.Range("A3").Value = "Shelly Safrai"
$newSheet.Range("A4").Value = "Yuval Melamed"

# The author also happened to leave the cursor on cell B29 of sheet "12"
# (navigated there while reviewing, outside its A1:A4 data range).
$sheet12 = $wb.Worksheets.Item("12")
$sheet12.Range("B29").Select() | Out-Null

# Finally, re-activate the new sheet so it stays the selected/visible tab,
# with its cursor left on the last data row (A4).
$newSheet.Activate() | Out-Null
$newSheet.Range("A4").Select() | Out-Null
